$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:ASHAPURMIN"
$ws.Range("C2").Value = "NSE:AARTIIND"
$ws.Range("D2").Value = "NSE:JUBLFOOD"
$ws.Range("E2").Value = "NSE:NATIONALUM"
$ws.Range("F2").Value = "NSE:BAJAJFINSV"
$ws.Range("B3").Value = "NSE:BAJFINANCE"
$ws.Range("C3").Value = "NSE:AKSHOPTFBR"
$ws.Range("D3").Value = "NSE:M&M"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "NSE:BAJFINANCE"
$ws.Range("B4").Value = "NSE:BAYERCROP"
$ws.Range("C4").Value = "NSE:ALANKIT"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "NSE:HINDPETRO"
$ws.Range("B5").Value = "NSE:BEDMUTHA"
$ws.Range("C5").Value = "NSE:APTUS"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "NSE:MAXHEALTH"
$ws.Range("B6").Value = "NSE:BIRLACORPN"
$ws.Range("C6").Value = "NSE:ARTNIRMAN"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("B7").Value = "NSE:CUMMINSIND"
$ws.Range("C7").Value = "NSE:ARVSMART"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("B8").Value = "NSE:CYIENT"
$ws.Range("C8").Value = "NSE:BANARISUG"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("B9").Value = "NSE:ERIS"
$ws.Range("C9").Value = "NSE:BANKINDIA"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("B10").Value = "NSE:ESTER"
$ws.Range("C10").Value = "NSE:BIOCON"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("B11").Value = "NSE:GIPCL"
$ws.Range("C11").Value = "NSE:BPL"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("B12").Value = "NSE:GOACARBON"
$ws.Range("C12").Value = "NSE:DEEPAKFERT"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("B13").Value = "NSE:GODFRYPHLP"
$ws.Range("C13").Value = "NSE:DYCL"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("B14").Value = "NSE:HINDCOMPOS"
$ws.Range("C14").Value = "NSE:ELGIRUBCO"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("B15").Value = "NSE:HINDPETRO"
$ws.Range("C15").Value = "NSE:EMMBI"
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""
$ws.Range("B16").Value = "NSE:HIRECT"
$ws.Range("C16").Value = "NSE:FIBERWEB"
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""
$ws.Range("B17").Value = "NSE:INDIAMART"
$ws.Range("C17").Value = "NSE:GENSOL"
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("B18").Value = "NSE:INDIANHUME"
$ws.Range("C18").Value = "NSE:GTL"
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("B19").Value = "NSE:JKLAKSHMI"
$ws.Range("C19").Value = "NSE:HPAL"
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("B20").Value = "NSE:JUBLINGREA"
$ws.Range("C20").Value = "NSE:INDOWIND"
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("B21").Value = "NSE:LATENTVIEW"
$ws.Range("C21").Value = "NSE:INDUSTOWER"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = ""
$ws.Range("B22").Value = "NSE:LGHL"
$ws.Range("C22").Value = "NSE:NUCLEUS"
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = ""
$ws.Range("B23").Value = "NSE:MAXHEALTH"
$ws.Range("C23").Value = "NSE:OBCL"
$ws.Range("D23").Value = ""
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = ""
$ws.Range("B24").Value = "NSE:MON100"
$ws.Range("C24").Value = "NSE:PAYTM"
$ws.Range("D24").Value = ""
$ws.Range("E24").Value = ""
$ws.Range("F24").Value = ""
$ws.Range("B25").Value = "NSE:NKIND"
$ws.Range("C25").Value = "NSE:PRAKASHSTL"
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""
$ws.Range("B26").Value = "NSE:PKTEA"
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = ""
$ws.Range("E26").Value = ""
$ws.Range("F26").Value = ""
$ws.Range("B27").Value = "NSE:RHIM"
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = ""
$ws.Range("F27").Value = ""

# Remove the now-unused rows 28-33 (sheet shrinks from F33 to F27)
$ws.Range("A28:F33").Delete()
